# Applies: "Added first Conditional Token to code generator routine"
#
# 1. tokens sheet: add a new [option_ids] section (header + count + one
#    token row) defining the first conditional token, %{(101)}%, bound to
#    the new `pterm_enable` option.
# 2. comp_header sheet: prefix the P-Term controller declaration block
#    (8 lines) with the new %{(101)}% conditional token so that block is
#    now only emitted when option 101 (pterm_enable) is active.
# 3. View-state bookkeeping: make "tokens" the active sheet/tab again
#    (matches activeTab/tabSelected move away from library_header), and
#    update a couple of stored cell selections to match where the author
#    ended up after editing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) tokens sheet - new [option_ids] section
# ---------------------------------------------------------------------
$tokens = $wb.Worksheets.Item("tokens")

# Row 2's existing count formula is rewritten by Excel from the
# "COUNTA(A:A)-2" idiom to an explicit range once a second counted
# section exists further down the sheet.
$tokens.Range("C2").Formula = "=COUNTA(A2:A20)"

# Clone the formatting of the existing "[tokens]" section (rows 1-3) down
# onto the new section (rows 21-23) before filling in values, so the new
# section header/count/data rows keep the same look (bold white-on-navy
# header, etc.).
$tokens.Range("A1:C1").Copy() | Out-Null
$tokens.Range("A21:C21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$tokens.Range("A2:C2").Copy() | Out-Null
$tokens.Range("A22:C22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$tokens.Range("A3:C3").Copy() | Out-Null
$tokens.Range("A23:C23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# New section header
$tokens.Range("A21").Value = "[option_ids]"

# New section count row
$tokens.Range("A22").Value = "count"
$tokens.Range("B22").Value = "'="
$tokens.Range("C22").Formula = "=COUNTA(A23:A30)"

# First (and so far only) conditional-token / option id entry
$tokens.Range("A23").Value = 0
$tokens.Range("B23").Value = "'="
$tokens.Range("C23").Value = "%{(101)}%;pterm_enable"

# ---------------------------------------------------------------------
# 2) comp_header sheet - gate the P-Term controller declaration behind
#    the new %{(101)}% conditional token
# ---------------------------------------------------------------------
$compHeader = $wb.Worksheets.Item("comp_header")

$compHeader.Range("C88").Value = "%{(101)}%// Calls the %FILENAME_PATTERN% P-Term controller during measurements of plant transfer functions"
$compHeader.Range("C89").Value = "%{(101)}%// THIS CONTROLLER IS USED FOR MEASUREMENTS OF THE PLANT TRANSFER FUNCTION ONLY."
$compHeader.Range("C90").Value = "%{(101)}%// THIS LOOP IS BY DEFAULT UNSTABLE AND ONLY WORKS UNDER STABLE TEST CONDITIONS"
$compHeader.Range("C91").Value = "%{(101)}%// DO NOT USE THIS CONTROLLER TYPE FOR NORMAL OPERATION"
$compHeader.Range("C92").Value = "%{(101)}%extern void %FILENAME_PATTERN%_PTermUpdate( // Calls the P-Term controller (Assembly)"
$compHeader.Range("C93").Value = "%{(101)}%%IDENT%%IDENT%volatile %STRUCTURE_LABEL%* controller // Pointer to nPnZ data type object"
$compHeader.Range("C94").Value = "%{(101)}%%IDENT%);"
$compHeader.Range("C95").Value = "%{(101)}%%EMPTY%"

# ---------------------------------------------------------------------
# 3) View-state: author ended up with the tokens sheet selected/active
#    again, having left a particular cell selected on comp_header and
#    comp_source_head scrolled further down.
# ---------------------------------------------------------------------
$compHeader.Range("C96").Select() | Out-Null

$compSourceHead = $wb.Worksheets.Item("comp_source_head")
$compSourceHead.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 49

$tokens.Activate() | Out-Null
$tokens.Range("C24").Select() | Out-Null
